# This script re-applies a re-sorted/re-numbered export of the "Artfynd" sheet.
# A handful of rows had their entire contents (all columns A:AY) reshuffled
# among each other (their row numbers / formatting stay fixed, only the
# data moved). We implement this as a set of whole-row array swaps /
# rotations using full-row ranges (A:AY) read and written back via the
# Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($rowNum) {
    return $ws.Range("A" + $rowNum + ":AY" + $rowNum).Value()
}

function Set-RowValues($rowNum, $values) {
    $ws.Range("A" + $rowNum + ":AY" + $rowNum).Value = $values
}

function Swap-Rows($r1, $r2) {
    $v1 = Get-RowValues $r1
    $v2 = Get-RowValues $r2
    Set-RowValues $r1 $v2
    Set-RowValues $r2 $v1
}

# Apply a cyclic rotation of row contents where the content that was in
# $rows[0] ends up in $rows[1], the content that was in $rows[1] ends up
# in $rows[2], ..., and the content that was in the last row ends up in
# $rows[0].
# All source values are read up-front (before any writes happen) so that
# writes to one row in the group can never clobber a value that still
# needs to be read from that row.
function Rotate-Rows($rows) {
    $n = $rows.Length
    $saved = @()
    for ($i = 0; $i -lt $n; $i++) {
        $saved += , (Get-RowValues $rows[$i])
    }
    for ($i = 0; $i -lt $n; $i++) {
        $dst = $rows[($i + 1) % $n]
        Set-RowValues $dst $saved[$i]
    }
}

# Rows 25 <-> 26 (simple swap)
Swap-Rows 25 26

# Rows 32 <-> 33 (simple swap)
Swap-Rows 32 33

# Rows 45 -> 46 -> 47 -> 48 -> 45 (4-way cycle)
Rotate-Rows @(45, 46, 47, 48)

# Rows 49 <-> 50 (simple swap)
Swap-Rows 49 50

# Rows 61 -> 62 -> 63 -> 61 (3-way cycle)
Rotate-Rows @(61, 62, 63)

# Rows 70 -> 71 -> 72 -> 73 -> 70 (4-way cycle)
Rotate-Rows @(70, 71, 72, 73)

# Rows 82 <-> 83 (simple swap)
Swap-Rows 82 83
